$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.955131053924561
$ws.Range("B1").Value = 1.963433980941772
$ws.Range("C1").Value = 1.867145299911499
$ws.Range("D1").Value = 2.554807186126709
$ws.Range("E1").Value = 5.070760726928711
